# Adds RTs (reaction-time-style summary columns) for all trial types.
# New columns CU:DG (13 cols) appended after existing CT, for header row 1
# and data rows 2:47. Also updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1), columns CU:DG ----
$hdr = New-Object 'object[,]' 1,13
$hdr[0,0] = "LG"
$hdr[0,1] = "SG"
$hdr[0,2] = "N"
$hdr[0,3] = "SL"
$hdr[0,4] = "LL"
$hdr[0,5] = "V_beta"
$hdr[0,6] = "LG-N"
$hdr[0,7] = "SG-N"
$hdr[0,8] = "LL-N"
$hdr[0,9] = "SL_N"
$hdr[0,10] = "rel_rew"
$hdr[0,11] = "LG-SG"
$hdr[0,12] = "LL_SL"
$ws.Range("CU1:DG1").Value = $hdr

# ---- Data rows 2:47, columns CU:DG ----
$dat = New-Object 'object[,]' 46,13
$dat[0,0] = 0.305434638867154
$dat[0,1] = 0.32050712750060401
$dat[0,2] = 0.332187023188453
$dat[0,3] = 0.362249740224797
$dat[0,4] = 0.30837125348625699
$dat[0,5] = -0.0085370806711060794
$dat[0,6] = -0.026752384321298402
$dat[0,7] = -0.0116798956878483
$dat[0,8] = -0.0238157697021961
$dat[0,9] = 0.030062717036344101
$dat[0,10] = -0.0029366146191023199
$dat[0,11] = -0.0150724886334501
$dat[0,12] = -0.053878486738540197
$dat[1,0] = 0.37713703014014699
$dat[1,1] = 0.34409770040656401
$dat[1,2] = 0.32311334281985099
$dat[1,3] = 0.389765408923267
$dat[1,4] = 0.30858965820516399
$dat[1,5] = -0.00061688701992225795
$dat[1,6] = 0.0540236873202957
$dat[1,7] = 0.0209843575867125
$dat[1,8] = -0.014523684614687201
$dat[1,9] = 0.066652066103415494
$dat[1,10] = 0.068547371934982906
$dat[1,11] = 0.033039329733583103
$dat[1,12] = -0.081175750718102693
$dat[2,0] = 0.34347991854883703
$dat[2,1] = 0.33897643524687698
$dat[2,2] = 0.32272362557705397
$dat[2,3] = 0.34278606483712698
$dat[2,4] = 0.302930305479094
$dat[2,5] = -0.0024563787987322799
$dat[2,6] = 0.0207562929717823
$dat[2,7] = 0.016252809669822399
$dat[2,8] = -0.0197933200979605
$dat[2,9] = 0.020062439260072999
$dat[2,10] = 0.040549613069742897
$dat[2,11] = 0.0045034833019599301
$dat[2,12] = -0.039855759358033503
$dat[3,0] = 0.271466309786774
$dat[3,1] = 0.27853118465282001
$dat[3,2] = 0.294866868120152
$dat[3,3] = 0.29255829763133001
$dat[3,4] = 0.271871120959986
$dat[3,5] = -0.0052963112164953899
$dat[3,6] = -0.023400558333378198
$dat[3,7] = -0.016335683467332201
$dat[3,8] = -0.022995747160166499
$dat[3,9] = -0.0023085704888217099
$dat[3,10] = -0.00040481117321178301
$dat[3,11] = -0.0070648748660460097
$dat[3,12] = -0.0206871766713447
$dat[4,0] = 0.27033338195178602
$dat[4,1] = 0.27706846105866101
$dat[4,2] = 0.28496288182213902
$dat[4,3] = 0.27218974754214198
$dat[4,4] = 0.277509949752129
$dat[4,5] = -0.0016783792026607001
$dat[4,6] = -0.0146294998703524
$dat[4,7] = -0.0078944207634776796
$dat[4,8] = -0.0074529320700094104
$dat[4,9] = -0.012773134279996099
$dat[4,10] = -0.0071765678003430297
$dat[4,11] = -0.0067350791068747596
$dat[4,12] = 0.0053202022099867402
$dat[5,0] = 0.28031748597277301
$dat[5,1] = 0.31086552736814999
$dat[5,2] = 0.295282690494786
$dat[5,3] = 0.318008447065949
$dat[5,4] = 0.29928748914971898
$dat[5,5] = -0.0043021003699063101
$dat[5,6] = -0.0149652045220136
$dat[5,7] = 0.0155828368733637
$dat[5,8] = 0.0040047986549325201
$dat[5,9] = 0.022725756571162398
$dat[5,10] = -0.018970003176946099
$dat[5,11] = -0.030548041395377298
$dat[5,12] = -0.018720957916229901
$dat[6,0] = 0.260252866370137
$dat[6,1] = 0.29065276565961501
$dat[6,2] = 0.302461104583926
$dat[6,3] = 0.28393655980471499
$dat[6,4] = 0.27895491209346801
$dat[6,5] = -0.00722114126464084
$dat[6,6] = -0.042208238213788697
$dat[6,7] = -0.0118083389243111
$dat[6,8] = -0.023506192490458402
$dat[6,9] = -0.018524544779211199
$dat[6,10] = -0.018702045723330198
$dat[6,11] = -0.030399899289477599
$dat[6,12] = -0.0049816477112471997
$dat[7,0] = 0.33834768389351599
$dat[7,1] = 0.31042992940638198
$dat[7,2] = 0.36112862802110601
$dat[7,3] = 0.33093400194775302
$dat[7,4] = 0.30117361911106799
$dat[7,5] = -0.0060413272419412101
$dat[7,6] = -0.022780944127589398
$dat[7,7] = -0.050698698614723897
$dat[7,8] = -0.059955008910037501
$dat[7,9] = -0.030194626073352899
$dat[7,10] = 0.037174064782448099
$dat[7,11] = 0.027917754487134499
$dat[7,12] = -0.029760382836684501
$dat[8,0] = 0.270464802962123
$dat[8,1] = 0.27610288557389101
$dat[8,2] = 0.28752696858191401
$dat[8,3] = 0.29879739821717499
$dat[8,4] = 0.26873461768263901
$dat[8,5] = -0.0051110985475263798
$dat[8,6] = -0.0170621656197909
$dat[8,7] = -0.011424083008023399
$dat[8,8] = -0.018792350899275299
$dat[8,9] = 0.011270429635260299
$dat[8,10] = 0.0017301852794844299
$dat[8,11] = -0.0056380826117674501
$dat[8,12] = -0.030062780534535701
$dat[9,0] = 0.28478656522929602
$dat[9,1] = 0.322761286923196
$dat[9,2] = 0.29669485258636902
$dat[9,3] = 0.33304927079007002
$dat[9,4] = 0.28506881650537202
$dat[9,5] = -0.0078206785297620096
$dat[9,6] = -0.011908287357073199
$dat[9,7] = 0.026066434336826198
$dat[9,8] = -0.011626036080997401
$dat[9,9] = 0.036354418203700299
$dat[9,10] = -0.00028225127607584
$dat[9,11] = -0.037974721693899399
$dat[9,12] = -0.047980454284697702
$dat[10,0] = 0.24546278372872599
$dat[10,1] = 0.25314377364702501
$dat[10,2] = 0.25403052347246502
$dat[10,3] = 0.25497810076922101
$dat[10,4] = 0.249503193655982
$dat[10,5] = -0.0018750690422685899
$dat[10,6] = -0.0085677397437393596
$dat[10,7] = -0.00088674982544034698
$dat[10,8] = -0.0045273298164829603
$dat[10,9] = 0.00094757729675620697
$dat[10,10] = -0.0040404099272564001
$dat[10,11] = -0.0076809899182990098
$dat[10,12] = -0.0054749071132391596
$dat[11,0] = 0.29833598533878097
$dat[11,1] = 0.34576657786965298
$dat[11,2] = 0.315847488935105
$dat[11,3] = 0.33085254079196602
$dat[11,4] = 0.298982144508045
$dat[11,5] = -0.0081198454884412793
$dat[11,6] = -0.017511503596324401
$dat[11,7] = 0.029919088934548199
$dat[11,8] = -0.016865344427060301
$dat[11,9] = 0.015005051856860499
$dat[11,10] = -0.00064615916926413699
$dat[11,11] = -0.047430592530872603
$dat[11,12] = -0.031870396283920799
$dat[12,0] = 0.55924635590054095
$dat[12,1] = 0.52810588863212604
$dat[12,2] = 0.48599444975843598
$dat[12,3] = 0.57096491596894305
$dat[12,4] = 0.478740440797992
$dat[12,5] = 0.00035099209136602
$dat[12,6] = 0.073251906142104403
$dat[12,7] = 0.042111438873689601
$dat[12,8] = -0.0072540089604444802
$dat[12,9] = 0.084970466210506801
$dat[12,10] = 0.080505915102548897
$dat[12,11] = 0.031140467268414702
$dat[12,12] = -0.092224475170951295
$dat[13,0] = 0.27883449981891301
$dat[13,1] = 0.25725000037527901
$dat[13,2] = 0.29304862715798602
$dat[13,3] = 0.34862338516904801
$dat[13,4] = 0.305491818497102
$dat[13,5] = -0.00166557165916207
$dat[13,6] = -0.0142141273390734
$dat[13,7] = -0.035798626782707198
$dat[13,8] = 0.0124431913391163
$dat[13,9] = 0.055574758011061903
$dat[13,10] = -0.026657318678189702
$dat[13,11] = 0.0215844994436338
$dat[13,12] = -0.0431315666719456
$dat[14,0] = 0.26968774611941598
$dat[14,1] = 0.29450151343735298
$dat[14,2] = 0.28823237660981199
$dat[14,3] = 0.268033032094535
$dat[14,4] = 0.28180007770151799
$dat[14,5] = -0.00257311793640319
$dat[14,6] = -0.018544630490396199
$dat[14,7] = 0.0062691368275409297
$dat[14,8] = -0.0064322989082938797
$dat[14,9] = -0.0201993445152766
$dat[14,10] = -0.012112331582102299
$dat[14,11] = -0.024813767317937101
$dat[14,12] = 0.013767045606982701
$dat[15,0] = 0.29193674297130201
$dat[15,1] = 0.29208420826762399
$dat[15,2] = 0.32701431243913198
$dat[15,3] = 0.29602516900922599
$dat[15,4] = 0.33810273605922703
$dat[15,5] = 0.0012814968504244299
$dat[15,6] = -0.035077569467830402
$dat[15,7] = -0.034930104171507902
$dat[15,8] = 0.0110884236200945
$dat[15,9] = -0.030989143429906101
$dat[15,10] = -0.046165993087925003
$dat[15,11] = -0.000147465296322479
$dat[15,12] = 0.042077567050000597
$dat[16,0] = 0.21253948134835801
$dat[16,1] = 0.227537484490312
$dat[16,2] = 0.25585869397036698
$dat[16,3] = 0.21015316783450499
$dat[16,4] = 0.22880022309254799
$dat[16,5] = -0.0047663308131242199
$dat[16,6] = -0.043319212622009197
$dat[16,7] = -0.0283212094800546
$dat[16,8] = -0.027058470877818701
$dat[16,9] = -0.045705526135861797
$dat[16,10] = -0.016260741744190399
$dat[16,11] = -0.0149980031419545
$dat[16,12] = 0.0186470552580431
$dat[17,0] = 0.30696640536189002
$dat[17,1] = 0.29056982032488998
$dat[17,2] = 0.30336897366214499
$dat[17,3] = 0.30574245750904
$dat[17,4] = 0.30899363989010398
$dat[17,5] = 0.00206213323898347
$dat[17,6] = 0.0035974316997453501
$dat[17,7] = -0.012799153337255099
$dat[17,8] = 0.0056246662279590904
$dat[17,9] = 0.0023734838468953902
$dat[17,10] = -0.0020272345282137299
$dat[17,11] = 0.016396585037000401
$dat[17,12] = 0.0032511823810636902
$dat[18,0] = 0.23650924318644601
$dat[18,1] = 0.23949521486974801
$dat[18,2] = 0.26371218028361898
$dat[18,3] = 0.255010404034692
$dat[18,4] = 0.24415008858522899
$dat[18,5] = -0.0043293797091662001
$dat[18,6] = -0.027202937097172201
$dat[18,7] = -0.024216965413870601
$dat[18,8] = -0.019562091698389801
$dat[18,9] = -0.0087017762489267608
$dat[18,10] = -0.0076408453987823997
$dat[18,11] = -0.0029859716833016102
$dat[18,12] = -0.010860315449463
$dat[19,0] = 0.31192951760749499
$dat[19,1] = 0.316161317030491
$dat[19,2] = 0.30807536066276903
$dat[19,3] = 0.30272618021990599
$dat[19,4] = 0.31385727528686402
$dat[19,5] = 0.0011810976580558701
$dat[19,6] = 0.0038541569447261199
$dat[19,7] = 0.0080859563677222399
$dat[19,8] = 0.0057819146240944896
$dat[19,9] = -0.0053491804428631396
$dat[19,10] = -0.00192775767936836
$dat[19,11] = -0.00423179942299611
$dat[19,12] = 0.0111310950669576
$dat[20,0] = 0.276945968158543
$dat[20,1] = 0.27892927033826698
$dat[20,2] = 0.27071501815225901
$dat[20,3] = 0.27823677496053201
$dat[20,4] = 0.29237708583241301
$dat[20,5] = 0.0028607161698996002
$dat[20,6] = 0.0062309500062838197
$dat[20,7] = 0.0082142521860077908
$dat[20,8] = 0.021662067680153899
$dat[20,9] = 0.0075217568082734899
$dat[20,10] = -0.0154311176738701
$dat[20,11] = -0.0019833021797239698
$dat[20,12] = 0.0141403108718805
$dat[21,0] = 0.31653699734078999
$dat[21,1] = 0.29115420039033701
$dat[21,2] = 0.28149277151260299
$dat[21,3] = 0.31675721331794099
$dat[21,4] = 0.31341367367713202
$dat[21,5] = 0.0063574560930256801
$dat[21,6] = 0.035044225828187302
$dat[21,7] = 0.0096614288777345791
$dat[21,8] = 0.031920902164529197
$dat[21,9] = 0.035264441805338699
$dat[21,10] = 0.0031233236636580802
$dat[21,11] = 0.025382796950452698
$dat[21,12] = -0.0033435396408094602
$dat[22,0] = 0.29694211028981898
$dat[22,1] = 0.239452081179479
$dat[22,2] = 0.219973833620315
$dat[22,3] = 0.22913358637015299
$dat[22,4] = 0.25394911135663201
$dat[22,5] = 0.0138035077501886
$dat[22,6] = 0.076968276669504093
$dat[22,7] = 0.019478247559163699
$dat[22,8] = 0.033975277736317297
$dat[22,9] = 0.0091597527498379298
$dat[22,10] = 0.042992998933186699
$dat[22,11] = 0.057490029110340402
$dat[22,12] = 0.0248155249864794
$dat[23,0] = 0.33442128088790901
$dat[23,1] = 0.31636660569347402
$dat[23,2] = 0.26959392149001299
$dat[23,3] = 0.31192002794705298
$dat[23,4] = 0.32269746973179197
$dat[23,5] = 0.0104830731870606
$dat[23,6] = 0.064827359397895606
$dat[23,7] = 0.046772684203460799
$dat[23,8] = 0.053103548241779201
$dat[23,9] = 0.0423261064570397
$dat[23,10] = 0.0117238111561164
$dat[23,11] = 0.018054675194434801
$dat[23,12] = 0.010777441784739401
$dat[24,0] = 0.35309921624138901
$dat[24,1] = 0.38635392149444597
$dat[24,2] = 0.370830671745352
$dat[24,3] = 0.38956117839552401
$dat[24,4] = 0.38565342366928201
$dat[24,5] = -0.0028622259685237698
$dat[24,6] = -0.017731455503962899
$dat[24,7] = 0.015523249749094199
$dat[24,8] = 0.0148227519239299
$dat[24,9] = 0.018730506650172101
$dat[24,10] = -0.032554207427892799
$dat[24,11] = -0.033254705253057099
$dat[24,12] = -0.0039077547262422697
$dat[25,0] = 0.309480648895259
$dat[25,1] = 0.32675396429840398
$dat[25,2] = 0.308012943074572
$dat[25,3] = 0.32442380778957097
$dat[25,4] = 0.30053509975550602
$dat[25,5] = -0.00336944006682772
$dat[25,6] = 0.00146770582068711
$dat[25,7] = 0.018741021223831901
$dat[25,8] = -0.0074778433190658602
$dat[25,9] = 0.0164108647149987
$dat[25,10] = 0.0089455491397529806
$dat[25,11] = -0.0172733154031448
$dat[25,12] = -0.0238887080340646
$dat[26,0] = 0.31061711788061003
$dat[26,1] = 0.327592697198269
$dat[26,2] = 0.31018091001897102
$dat[26,3] = 0.25669941294472598
$dat[26,4] = 0.31103491142857798
$dat[26,5] = 0.0027607234598170599
$dat[26,6] = 0.000436207861639559
$dat[26,7] = 0.017411787179298699
$dat[26,8] = 0.00085400140960700799
$dat[26,9] = -0.053481497074244502
$dat[26,10] = -0.00041779354796744802
$dat[26,11] = -0.016975579317659101
$dat[26,12] = 0.054335498483851503
$dat[27,0] = 0.24043969999911399
$dat[27,1] = 0.26800750000256801
$dat[27,2] = 0.27326109999921699
$dat[27,3] = 0.28887400000166902
$dat[27,4] = 0.258038300002226
$dat[27,5] = -0.0076034071428564103
$dat[27,6] = -0.032821400000102502
$dat[27,7] = -0.0052535999966494204
$dat[27,8] = -0.0152227999969909
$dat[27,9] = 0.015612900002451999
$dat[27,10] = -0.0175986000031116
$dat[27,11] = -0.027567800003453099
$dat[27,12] = -0.030835699999443002
$dat[28,0] = 0.30175958454492502
$dat[28,1] = 0.276141779846511
$dat[28,2] = 0.29592167254304502
$dat[28,3] = 0.285913090308895
$dat[28,4] = 0.27923204848775601
$dat[28,5] = 0.00057750363027611704
$dat[28,6] = 0.0058379120018798797
$dat[28,7] = -0.0197798926965333
$dat[28,8] = -0.016689624055288701
$dat[28,9] = -0.0100085822341497
$dat[28,10] = 0.0225275360571686
$dat[28,11] = 0.0256178046984132
$dat[28,12] = -0.0066810418211389298
$dat[29,0] = 0.24958635913208099
$dat[29,1] = 0.25142989447340303
$dat[29,2] = 0.25946936244145002
$dat[29,3] = 0.237753561581484
$dat[29,4] = 0.25313970015849901
$dat[29,5] = -0.00019071873975918101
$dat[29,6] = -0.0098830033093690803
$dat[29,7] = -0.00803946796804666
$dat[29,8] = -0.0063296622829511701
$dat[29,9] = -0.021715800859965301
$dat[29,10] = -0.0035533410264179102
$dat[29,11] = -0.0018435353413224201
$dat[29,12] = 0.015386138577014199
$dat[30,0] = 0.40755706885829501
$dat[30,1] = 0.41934548062272298
$dat[30,2] = 0.43431985168717802
$dat[30,3] = 0.373613651841878
$dat[30,4] = 0.32261203974485397
$dat[30,5] = -0.0143757584737613
$dat[30,6] = -0.026762782828882298
$dat[30,7] = -0.0149743710644543
$dat[30,8] = -0.111707811942324
$dat[30,9] = -0.060706199845299097
$dat[30,10] = 0.084945029113441706
$dat[30,11] = -0.011788411764428
$dat[30,12] = -0.051001612097024897
$dat[31,0] = 0.29532931579160498
$dat[31,1] = 0.31495051832462101
$dat[31,2] = 0.35074371173686802
$dat[31,3] = 0.29348909223335701
$dat[31,4] = 0.28530619526281897
$dat[31,5] = -0.010618286565919
$dat[31,6] = -0.055414395945262998
$dat[31,7] = -0.035793193412246099
$dat[31,8] = -0.065437516474048593
$dat[31,9] = -0.057254619503510101
$dat[31,10] = 0.0100231205287855
$dat[31,11] = -0.019621202533016899
$dat[31,12] = -0.0081828969705384207
$dat[32,0] = 0.25569203210761698
$dat[32,1] = 0.27452464896487
$dat[32,2] = 0.31919694674434101
$dat[32,3] = 0.21977720604627299
$dat[32,4] = 0.27808915924106198
$dat[32,5] = -0.0046523832716047599
$dat[32,6] = -0.063504914636723697
$dat[32,7] = -0.044672297779470599
$dat[32,8] = -0.0411077875032788
$dat[32,9] = -0.099419740698067402
$dat[32,10] = -0.0223971271334448
$dat[32,11] = -0.018832616857253001
$dat[32,12] = 0.058311953194788602
$dat[33,0] = 0.27095711999572802
$dat[33,1] = 0.26093404786661201
$dat[33,2] = 0.32185583957470898
$dat[33,3] = 0.29167764727026202
$dat[33,4] = 0.31236949667800201
$dat[33,5] = -0.00211929578134524
$dat[33,6] = -0.0508987195789814
$dat[33,7] = -0.060921791708096799
$dat[33,8] = -0.0094863428967073508
$dat[33,9] = -0.0301781923044472
$dat[33,10] = -0.041412376682274002
$dat[33,11] = 0.0100230721291154
$dat[33,12] = 0.020691849407739899
$dat[34,0] = 0.32115699999849301
$dat[34,1] = 0.32137630000943301
$dat[34,2] = 0.31825130000652202
$dat[34,3] = 0.33935379999456899
$dat[34,4] = 0.32835140000679502
$dat[34,5] = 0.00012743571382348899
$dat[34,6] = 0.00290569999197032
$dat[34,7] = 0.00312500000291038
$dat[34,8] = 0.0101001000002725
$dat[34,9] = 0.021102499988046398
$dat[34,10] = -0.0071944000083021802
$dat[34,11] = -0.00021930001094005999
$dat[34,12] = -0.0110023999877739
$dat[35,0] = 0.330868997305515
$dat[35,1] = 0.34617498883744702
$dat[35,2] = 0.32600145443575401
$dat[35,3] = 0.32111776137025999
$dat[35,4] = 0.33702493194141397
$dat[35,5] = 0.0011780142439030499
$dat[35,6] = 0.0048675428697606497
$dat[35,7] = 0.020173534401692401
$dat[35,8] = 0.0110234775056596
$dat[35,9] = -0.0048836930654942903
$dat[35,10] = -0.0061559346358990297
$dat[35,11] = -0.0153059915319317
$dat[35,12] = 0.0159071705711539
$dat[36,0] = 0.26031169999623599
$dat[36,1] = 0.24865109997335799
$dat[36,2] = 0.36765630001900701
$dat[36,3] = 0.25837260001571799
$dat[36,4] = 0.27290470001753397
$dat[36,5] = -0.0125645357142535
$dat[36,6] = -0.10734460002277001
$dat[36,7] = -0.119005200045648
$dat[36,8] = -0.094751600001472897
$dat[36,9] = -0.10928370000328801
$dat[36,10] = -0.0125930000212974
$dat[36,11] = 0.011660600022878401
$dat[36,12] = 0.0145321000018157
$dat[37,0] = 0.29838525387458498
$dat[37,1] = 0.26985317224171002
$dat[37,2] = 0.28624831972410902
$dat[37,3] = 0.30872106592869297
$dat[37,4] = 0.28544971602968799
$dat[37,5] = 0.00118564729928038
$dat[37,6] = 0.012136934150476
$dat[37,7] = -0.0163951474823988
$dat[37,8] = -0.00079860369442030701
$dat[37,9] = 0.022472746204584799
$dat[37,10] = 0.012935537844896299
$dat[37,11] = 0.028532081632874899
$dat[37,12] = -0.0232713498990051
$dat[38,0] = 0.26389229999767799
$dat[38,1] = 0.25869560000137398
$dat[38,2] = 0.27253230000496798
$dat[38,3] = 0.27097800000046801
$dat[38,4] = 0.262517000002844
$dat[38,5] = -0.0015656857150523501
$dat[38,6] = -0.0086400000072899205
$dat[38,7] = -0.013836700003594099
$dat[38,8] = -0.010015300002123599
$dat[38,9] = -0.0015543000044999601
$dat[38,10] = 0.00137529999483376
$dat[38,11] = 0.00519669999630423
$dat[38,12] = -0.0084609999976237305
$dat[39,0] = 0.27478430001065102
$dat[39,1] = 0.30509630002779797
$dat[39,2] = 0.42112389998510402
$dat[39,3] = 0.43343699997058099
$dat[39,4] = 0.34199659997830101
$dat[39,5] = -0.024801378570762998
$dat[39,6] = -0.14633959997445301
$dat[39,7] = -0.11602759995730499
$dat[39,8] = -0.079127300006803097
$dat[39,9] = 0.0123130999854765
$dat[39,10] = -0.067212299967650296
$dat[39,11] = -0.030312000017147501
$dat[39,12] = -0.091440399992279695
$dat[40,0] = 0.25526109989732498
$dat[40,1] = 0.32960830000229102
$dat[40,2] = 0.50983945000916697
$dat[40,3] = 0.31709309993311702
$dat[40,4] = 0.32746299984864802
$dat[40,5] = -0.035780864318699701
$dat[40,6] = -0.25457835011184199
$dat[40,7] = -0.18023115000687501
$dat[40,8] = -0.18237645016051801
$dat[40,9] = -0.19274635007604901
$dat[40,10] = -0.072201899951323797
$dat[40,11] = -0.074347200104966704
$dat[40,12] = 0.0103698999155312
$dat[41,0] = 0.25166589999571398
$dat[41,1] = 0.26142870000330698
$dat[41,2] = 0.26647029997548
$dat[41,3] = 0.25061879999702702
$dat[41,4] = 0.240791800024453
$dat[41,5] = -0.0042909071364972103
$dat[41,6] = -0.014804399979766399
$dat[41,7] = -0.0050415999721735699
$dat[41,8] = -0.025678499951027299
$dat[41,9] = -0.0158514999784529
$dat[41,10] = 0.010874099971260799
$dat[41,11] = -0.0097628000075928797
$dat[41,12] = -0.0098269999725744094
$dat[42,0] = 0.309154100017622
$dat[42,1] = 0.35540389991365301
$dat[42,2] = 0.44437369983643199
$dat[42,3] = 0.38728519994765498
$dat[42,4] = 0.32651179982349199
$dat[42,5] = -0.025721764275138899
$dat[42,6] = -0.13521959981880999
$dat[42,7] = -0.088969799922779202
$dat[42,8] = -0.11786190001294
$dat[42,9] = -0.057088499888777698
$dat[42,10] = -0.0173576998058706
$dat[42,11] = -0.046249799896031597
$dat[42,12] = -0.060773400124162401
$dat[43,0] = 0.26531110005453201
$dat[43,1] = 0.30091230000834901
$dat[43,2] = 0.31645939999725597
$dat[43,3] = 0.37248519994318402
$dat[43,4] = 0.26834149996284301
$dat[43,5] = -0.017072221422235301
$dat[43,6] = -0.0511482999427244
$dat[43,7] = -0.015547099988907499
$dat[43,8] = -0.0481179000344127
$dat[43,9] = 0.056025799945928101
$dat[43,10] = -0.0030303999083116599
$dat[43,11] = -0.035601199953816803
$dat[43,12] = -0.10414369998034
$dat[44,0] = 0.29868469992652502
$dat[44,1] = 0.27217550016939601
$dat[44,2] = 0.299600899685174
$dat[44,3] = 0.30233040032908298
$dat[44,4] = 0.27748880023136702
$dat[44,5] = -0.0015257642537887699
$dat[44,6] = -0.00091619975864887205
$dat[44,7] = -0.027425399515777799
$dat[44,8] = -0.022112099453806801
$dat[44,9] = 0.0027295006439089701
$dat[44,10] = 0.021195899695158001
$dat[44,11] = 0.026509199757128898
$dat[44,12] = -0.024841600097715799
$dat[45,0] = 0.35358945000916697
$dat[45,1] = 0.45092039974406301
$dat[45,2] = 0.37348589999601201
$dat[45,3] = 0.44888090016320298
$dat[45,4] = 0.34719959995709299
$dat[45,5] = -0.017513928569054999
$dat[45,6] = -0.019896449986845199
$dat[45,7] = 0.077434499748051097
$dat[45,8] = -0.026286300038918801
$dat[45,9] = 0.075395000167191001
$dat[45,10] = 0.0063898500520735901
$dat[45,11] = -0.097330949734896394
$dat[45,12] = -0.10168130020610899
$ws.Range("CU2:DG47").Value = $dat

# ---- Restore active selection on the (now wider) sheet ----
$ws.Range("DH26").Select() | Out-Null
